# Merge the split "<id>...</id>" runs for p033v_1 and p033v_2 into a
# single run each (keeping the first run's formatting/identity), i.e.
# turn:
#   [<id>](Courier New/7f6000) [p033v_N](black) [</id>](Courier New/7f6000)
# into:
#   [<id>p033v_N</id>](Courier New/7f6000)
#
# Note: "<id>fig_p033v_1</id>" elsewhere in the document must stay as-is.

$d = $word.ActiveDocument

function Merge-IdRuns([string]$needle) {
    $rng = $d.Content
    $found = $rng.Find.Execute($needle, $false, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
    if (-not $found) {
        return
    }

    $matchStart = $rng.Start
    $matchEnd = $rng.End

    # The match always starts with the literal "<id>" run boundary.
    $afterOpenTag = $matchStart + 4

    # Grab (and remember) everything after "<id>" within the match, then
    # delete it so only the first run ("<id>") is left.
    $rest = $d.Range($afterOpenTag, $matchEnd)
    $restText = $rest.Text
    $rest.Delete()

    # Re-insert the remaining text right at the (now collapsed) boundary
    # so it is appended into the existing first run instead of spawning
    # a freshly-styled one.
    $insertionPoint = $d.Range($afterOpenTag, $afterOpenTag)
    $insertionPoint.InsertAfter($restText)
}

Merge-IdRuns "<id>p033v_1</id>"
Merge-IdRuns "<id>p033v_2</id>"
